$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Range("D38").Value = 254
$ws.Range("L38").Value = 1365.125
$ws.Range("N38").Value = -82.875

# Row 39
$ws.Range("D39").Value = 254
$ws.Range("L39").Value = 1365.125
$ws.Range("N39").Value = -34.875

# Row 40
$ws.Range("D40").Value = 254
$ws.Range("L40").Value = 1365.125
$ws.Range("N40").Value = -73.875

# Row 41
$ws.Range("D41").Value = 254
$ws.Range("L41").Value = 1365.125
$ws.Range("N41").Value = -44.47499999999991

# Row 42
$ws.Range("B42").Value = 485.025
$ws.Range("D42").Value = 317
$ws.Range("H42").Value = 131
$ws.Range("J42").Value = 196
$ws.Range("L42").Value = 1129.025
$ws.Range("N42").Value = -97.97499999999991

# Row 43
$ws.Range("B43").Value = 485.025
$ws.Range("D43").Value = 317
$ws.Range("G43").Value = 1210.4
$ws.Range("H43").Value = 131
$ws.Range("J43").Value = 196
$ws.Range("L43").Value = 1129.025
$ws.Range("M43").Value = 1210.4
$ws.Range("N43").Value = -81.375

# Row 44
$ws.Range("B44").Value = 485.025
$ws.Range("D44").Value = 317
$ws.Range("G44").Value = 1164.5
$ws.Range("H44").Value = 131
$ws.Range("J44").Value = 196
$ws.Range("L44").Value = 1129.025
$ws.Range("M44").Value = 1164.5
$ws.Range("N44").Value = -35.47499999999991

# Row 45
$ws.Range("B45").Value = 485.025
$ws.Range("D45").Value = 317
$ws.Range("G45").Value = 1104.9
$ws.Range("H45").Value = 131
$ws.Range("J45").Value = 196
$ws.Range("L45").Value = 1129.025
$ws.Range("M45").Value = 1104.9
$ws.Range("N45").Value = 24.125

# Row 46
$ws.Range("B46").Value = 429.8
$ws.Range("D46").Value = 137
$ws.Range("G46").Value = 1033.6
$ws.Range("H46").Value = 152
$ws.Range("J46").Value = 244
$ws.Range("L46").Value = 962.8
$ws.Range("M46").Value = 1033.6
$ws.Range("N46").Value = -70.79999999999995

# Row 47
$ws.Range("B47").Value = 429.8
$ws.Range("D47").Value = 137
$ws.Range("G47").Value = 959.9000000000001
$ws.Range("H47").Value = 152
$ws.Range("J47").Value = 244
$ws.Range("L47").Value = 962.8
$ws.Range("M47").Value = 959.9000000000001
$ws.Range("N47").Value = 2.899999999999864

# Row 48
$ws.Range("B48").Value = 429.8
$ws.Range("D48").Value = 137
$ws.Range("G48").Value = 993.7
$ws.Range("H48").Value = 152
$ws.Range("J48").Value = 244
$ws.Range("L48").Value = 962.8
$ws.Range("M48").Value = 993.7
$ws.Range("N48").Value = -30.90000000000009

# Row 49
$ws.Range("B49").Value = 429.8
$ws.Range("D49").Value = 137
$ws.Range("G49").Value = 1041.5
$ws.Range("H49").Value = 152
$ws.Range("J49").Value = 244
$ws.Range("L49").Value = 962.8
$ws.Range("M49").Value = 1041.5
$ws.Range("N49").Value = -78.70000000000005

# Row 50
$ws.Range("G50").Value = 1090.2
$ws.Range("M50").Value = 1090.2
$ws.Range("N50").Value = -1090.2
